$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "CompleteSVM"

# Update headers (C1:F1)
$ws.Range("C1").Value = "Dataset"
$ws.Range("D1").Value = "Kernel"
$ws.Range("E1").Value = "C"
$ws.Range("F1").Value = "Gamma"

# Update data rows 2-17
$ws.Range("C2").Value = "Complete"
$ws.Range("D2").Value = "rbf"
$ws.Range("E2").Value = 0.001
$ws.Range("F2").Value = "scale"
$ws.Range("G2").Value = 0.1915484268969772
$ws.Range("H2").Value = 0.1865203761755486
$ws.Range("I2").Value = 0.1839443742098609
$ws.Range("J2").Value = 0.1837444655281467
$ws.Range("K2").Value = 0.1898263027295285
$ws.Range("L2").Value = 0.1871167891080124

$ws.Range("C3").Value = "Complete"
$ws.Range("D3").Value = "rbf"
$ws.Range("E3").Value = 0.01
$ws.Range("F3").Value = "scale"
$ws.Range("G3").Value = 0.6808098560552988
$ws.Range("H3").Value = 0.6767517672393788
$ws.Range("I3").Value = 0.6793398127806943
$ws.Range("J3").Value = 0.666273698520587
$ws.Range("K3").Value = 0.670045632106979
$ws.Range("L3").Value = 0.6746441533405877

$ws.Range("C4").Value = "Complete"
$ws.Range("D4").Value = "rbf"
$ws.Range("E4").Value = 0.05
$ws.Range("F4").Value = "scale"
$ws.Range("G4").Value = 0.9538322171813265
$ws.Range("H4").Value = 0.9489665752304639
$ws.Range("I4").Value = 0.9581775081535457
$ws.Range("J4").Value = 0.9601877803619121
$ws.Range("K4").Value = 0.9643856934411672
$ws.Range("L4").Value = 0.9571099548736832

$ws.Range("C5").Value = "Complete"
$ws.Range("D5").Value = "rbf"
$ws.Range("E5").Value = 0.1
$ws.Range("F5").Value = "scale"
$ws.Range("G5").Value = 0.9718942018997925
$ws.Range("H5").Value = 0.9638326935920185
$ws.Range("I5").Value = 0.971893679494582
$ws.Range("J5").Value = 0.9749476844002911
$ws.Range("K5").Value = 0.9754344732906393
$ws.Range("L5").Value = 0.9716005465354648

$ws.Range("C6").Value = "Complete"
$ws.Range("D6").Value = "rbf"
$ws.Range("E6").Value = 0.5
$ws.Range("F6").Value = "scale"
$ws.Range("G6").Value = 0.9867942008844881
$ws.Range("H6").Value = 0.9839995261344443
$ws.Range("I6").Value = 0.9883033215669523
$ws.Range("J6").Value = 0.9886695954093215
$ws.Range("K6").Value = 0.9878862125037932
$ws.Range("L6").Value = 0.9871305712997998

$ws.Range("C7").Value = "Complete"
$ws.Range("D7").Value = "rbf"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = "scale"
$ws.Range("G7").Value = 0.990265410527494
$ws.Range("H7").Value = 0.9874977058375438
$ws.Range("I7").Value = 0.9891673775946608
$ws.Range("J7").Value = 0.9911226669678685
$ws.Range("K7").Value = 0.9888643044701797
$ws.Range("L7").Value = 0.9893834930795492

$ws.Range("C8").Value = "Complete"
$ws.Range("D8").Value = "rbf"
$ws.Range("E8").Value = 1.5
$ws.Range("F8").Value = "scale"
$ws.Range("G8").Value = 0.9902157034865668
$ws.Range("H8").Value = 0.9880919409788653
$ws.Range("I8").Value = 0.9884427426472988
$ws.Range("J8").Value = 0.9909375427808046
$ws.Range("K8").Value = 0.99123273568448
$ws.Range("L8").Value = 0.989784133115603

$ws.Range("C9").Value = "Complete"
$ws.Range("D9").Value = "rbf"
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = "scale"
$ws.Range("G9").Value = 0.9919474211050463
$ws.Range("H9").Value = 0.9889594355443792
$ws.Range("I9").Value = 0.9892269793895476
$ws.Range("J9").Value = 0.9897325931167074
$ws.Range("K9").Value = 0.992212723776809
$ws.Range("L9").Value = 0.9904158305864978

$ws.Range("C10").Value = "Complete"
$ws.Range("D10").Value = "rbf"
$ws.Range("E10").Value = 0.001
$ws.Range("F10").Value = "auto"
$ws.Range("G10").Value = 0.1915484268969772
$ws.Range("H10").Value = 0.1865203761755486
$ws.Range("I10").Value = 0.1839443742098609
$ws.Range("J10").Value = 0.1837444655281467
$ws.Range("K10").Value = 0.1898263027295285
$ws.Range("L10").Value = 0.1871167891080124

$ws.Range("C11").Value = "Complete"
$ws.Range("D11").Value = "rbf"
$ws.Range("E11").Value = 0.01
$ws.Range("F11").Value = "auto"
$ws.Range("G11").Value = 0.6808098560552988
$ws.Range("H11").Value = 0.6767517672393788
$ws.Range("I11").Value = 0.6793398127806943
$ws.Range("J11").Value = 0.666273698520587
$ws.Range("K11").Value = 0.670045632106979
$ws.Range("L11").Value = 0.6746441533405877

$ws.Range("C12").Value = "Complete"
$ws.Range("D12").Value = "rbf"
$ws.Range("E12").Value = 0.05
$ws.Range("F12").Value = "auto"
$ws.Range("G12").Value = 0.9538322171813265
$ws.Range("H12").Value = 0.9489665752304639
$ws.Range("I12").Value = 0.9581775081535457
$ws.Range("J12").Value = 0.9601877803619121
$ws.Range("K12").Value = 0.9643856934411672
$ws.Range("L12").Value = 0.9571099548736832

$ws.Range("C13").Value = "Complete"
$ws.Range("D13").Value = "rbf"
$ws.Range("E13").Value = 0.1
$ws.Range("F13").Value = "auto"
$ws.Range("G13").Value = 0.9718942018997925
$ws.Range("H13").Value = 0.9638326935920185
$ws.Range("I13").Value = 0.971893679494582
$ws.Range("J13").Value = 0.9749476844002911
$ws.Range("K13").Value = 0.9754344732906393
$ws.Range("L13").Value = 0.9716005465354648

$ws.Range("C14").Value = "Complete"
$ws.Range("D14").Value = "rbf"
$ws.Range("E14").Value = 0.5
$ws.Range("F14").Value = "auto"
$ws.Range("G14").Value = 0.9867942008844881
$ws.Range("H14").Value = 0.9839995261344443
$ws.Range("I14").Value = 0.9883033215669523
$ws.Range("J14").Value = 0.9886695954093215
$ws.Range("K14").Value = 0.9878862125037932
$ws.Range("L14").Value = 0.9871305712997998

$ws.Range("C15").Value = "Complete"
$ws.Range("D15").Value = "rbf"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = "auto"
$ws.Range("G15").Value = 0.990265410527494
$ws.Range("H15").Value = 0.9874977058375438
$ws.Range("I15").Value = 0.9891673775946608
$ws.Range("J15").Value = 0.9911226669678685
$ws.Range("K15").Value = 0.9888643044701797
$ws.Range("L15").Value = 0.9893834930795492

$ws.Range("C16").Value = "Complete"
$ws.Range("D16").Value = "rbf"
$ws.Range("E16").Value = 1.5
$ws.Range("F16").Value = "auto"
$ws.Range("G16").Value = 0.9902157034865668
$ws.Range("H16").Value = 0.9880919409788653
$ws.Range("I16").Value = 0.9884427426472988
$ws.Range("J16").Value = 0.9909375427808046
$ws.Range("K16").Value = 0.99123273568448
$ws.Range("L16").Value = 0.989784133115603

$ws.Range("C17").Value = "Complete"
$ws.Range("D17").Value = "rbf"
$ws.Range("E17").Value = 5
$ws.Range("F17").Value = "auto"
$ws.Range("G17").Value = 0.9919474211050463
$ws.Range("H17").Value = 0.9889594355443792
$ws.Range("I17").Value = 0.9892269793895476
$ws.Range("J17").Value = 0.9897325931167074
$ws.Range("K17").Value = 0.992212723776809
$ws.Range("L17").Value = 0.9904158305864978

